# Termíno da Aula 04
# Appends text to the last paragraph and adds two new bullet paragraphs
# matching the existing "numId=1" list style used by the preceding items.

$d = $word.ActiveDocument

# --- 1. Extend the last paragraph's existing run with more text -----------
$lastParIndex = $d.Paragraphs.Count
$lastPar = $d.Paragraphs($lastParIndex)
$lastParRange = $lastPar.Range

# Collapse to just before the paragraph mark so InsertAfter appends onto the
# existing run instead of replacing the whole paragraph.
$appendPoint = $d.Range($lastParRange.Start, $lastParRange.End - 1)
[void]$appendPoint.InsertAfter(' Cada Entry do ChangeTracker, como já dito, tem uma referência para o objeto em si, e ele se encontra na propriedade Entity. Ao adicionar um objeto no contexto, o objeto fica com o estado Added nas Entries() do ChangeTracker e, ao chamar o saveChanges(), um objeto com esse estado é salvo no banco de dados através de um INSERT. Após isso, o objeto muda para o estado Unchanged.')

# --- 2. Insert two new list-item paragraphs after it -----------------------
# Re-fetch the (now extended) last paragraph range and collapse to its end
# (right after the paragraph mark) as the insertion point for new OOXML.
$lastPar = $d.Paragraphs($d.Paragraphs.Count)
$insertPoint = $d.Range($lastPar.Range.End, $lastPar.Range.End)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$pPrXml = '<w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:ind w:left="720" w:hanging="360"/><w:jc w:val="both"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="none"/></w:rPr></w:pPr>'
$rPrXml = '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr>'

$para1Xml = "<w:p $wNs>$pPrXml<w:r>$rPrXml<w:t xml:space=`"preserve`">Ao removermos um objeto do nosso contexto, sua Entry fica com o state “Deleted”. Qualquer objeto que possua sua Entry neste estado será apagado ao chamarmos o SaveChanges(). Um objeto, ao ser apagado no banco, não será mais monitorado pelo ChangeTracker.</w:t></w:r></w:p>"
$para2Xml = "<w:p $wNs>$pPrXml<w:r>$rPrXml<w:t xml:space=`"preserve`">Com o metódo de contexto Entry(), contexto.Entry(objeto), podemos obter a Entry de um determinado objeto. Se adicionarmos um produto no contexto e o apargamos sem chamar o SaveChanges(), ele ficará primeiro como Added e, ao remover, ele sumirá da lista. Ainda podemos obter o State desse objeto, com o metódo Entry(). Podemos ver que o state deste objeto é “Detached”.</w:t></w:r></w:p>"

[void]$insertPoint.InsertXML($para1Xml + $para2Xml)

Write-Output ("Paragraphs after edit: " + $d.Paragraphs.Count)

